$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Hydrogen): update B3, clear D3 (keep the cell present but blank)
$ws.Range("B3").Value = 497400.1395342923
$ws.Range("D3").ClearContents()
$ws.Range("D4").Copy()
$ws.Range("D3").PasteSpecial(-4122)

# Row 4 (Methanol): C4 -> 0
$ws.Range("C4").Value = 0

# Row 5 (Ammonia): C5 -> 2282.4042014741
$ws.Range("C5").Value = 2282.4042014741

# Row 7: rename "Other" -> "Biogas", update D7
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 1458.302950746877

# New row 8: "Other" label (matching style of the A-column labels), D8 value
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "Other"

$ws.Range("B7:C7").Copy()
$ws.Range("B8:C8").PasteSpecial(-4122)

$ws.Range("D8").Value = 1575.378061651642
